# Update header labels (row 1) across the workbook's sheets so that Power BI
# can automatically treat the first row as a header when the numeric-looking
# labels are loaded (they need a non-numeric prefix).
#
# - Sheets "Potencia Acumulada - SIN (MW)", "Geracao Periodo Medio (MWMed)",
#   "Atendimento a Ponta(MW)", "Emissoes Totais (MtCO2eq)" and
#   "Custo Total (bilhões de R$)" get an "Ano " prefix on their year labels.
# - Sheet "Potencia Incremental - SIN(MW)" gets an "Intervalo " prefix on its
#   year / interval labels.

$wb = $excel.ActiveWorkbook

$anoSheets = @(
    "Potencia Acumulada - SIN (MW)",
    "Geracao Periodo Medio (MWMed)",
    "Atendimento a Ponta(MW)",
    "Emissoes Totais (MtCO2eq)"
)

foreach ($sheetName in $anoSheets) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("B1").Value = "Ano 2015"
    $ws.Range("C1").Value = "Ano 2030"
    $ws.Range("D1").Value = "Ano 2040"
    $ws.Range("E1").Value = "Ano 2050"
}

$wsIncremental = $wb.Worksheets.Item("Potencia Incremental - SIN(MW)")
$wsIncremental.Range("B1").Value = "Intervalo 2015"
$wsIncremental.Range("C1").Value = "Intervalo 2015-2030"
$wsIncremental.Range("D1").Value = "Intervalo 2031-2040"
$wsIncremental.Range("E1").Value = "Intervalo 2041-2050"

$wsCusto = $wb.Worksheets.Item("Custo Total (bilhões de R$)")
$wsCusto.Range("B1").Value = "Ano 2015"
